$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'28.426.96"
$ws.Range('E2').Value = "'  +0.00%  "
$ws.Range('D3').Value = "'1.553.82"
$ws.Range('E3').Value = "'  -1.87%  "
$ws.Range('E4').Value = "'  -0.15%  "
$ws.Range('D5').Value = "'210.53"
$ws.Range('E5').Value = "'  -1.31%  "
$ws.Range('E6').Value = "'  -1.78%  "
$ws.Range('E7').Value = "'  -0.16%  "
$ws.Range('E8').Value = "'  +0.69%  "
$ws.Range('D9').Value = "'0.243"
$ws.Range('E9').Value = "'  -1.91%  "
$ws.Range('E10').Value = "'  -1.17%  "
$ws.Range('D11').Value = "'0.0890"
$ws.Range('E11').Value = "'  -0.43%  "
$ws.Range('D12').Value = "'1.775.80"
$ws.Range('E12').Value = "'  -1.92%  "
$ws.Range('D13').Value = "'1.554.86"
$ws.Range('E13').Value = "'  -1.85%  "
$ws.Range('D14').Value = "'28.432.64"
$ws.Range('E14').Value = "'  +0.00%  "
$ws.Range('E15').Value = "'  -2.07%  "
$ws.Range('E16').Value = "'  -1.59%  "
$ws.Range('E17').Value = "'  -1.56%  "
$ws.Range('D18').Value = "'228.88"
$ws.Range('E18').Value = "'  -0.49%  "
$ws.Range('E19').Value = "'  -1.72%  "
$ws.Range('E20').Value = "'  -2.24%  "
$ws.Range('E21').Value = "'  -0.10%  "
$ws.Range('D22').Value = "'3.88"
$ws.Range('E22').Value = "'  -0.61%  "
$ws.Range('D23').Value = "'8.91"
$ws.Range('D24').Value = "'2.02"
$ws.Range('E24').Value = "'  -2.02%  "
$ws.Range('D25').Value = "'151.41"
$ws.Range('E25').Value = "'  -0.30%  "
$ws.Range('E26').Value = "'  -2.01%  "
$ws.Range('E27').Value = "'  -1.28%  "
$ws.Range('E28').Value = "'  -0.13%  "
$ws.Range('E29').Value = "'  -3.06%  "
$ws.Range('E30').Value = "'  -2.94%  "
$ws.Range('E31').Value = "'  -4.43%  "
$ws.Range('E32').Value = "'  -1.49%  "
$ws.Range('D33').Value = "'1.384.11"
$ws.Range('E33').Value = "'  -0.84%  "
$ws.Range('E34').Value = "'  -3.16%  "
$ws.Range('E35').Value = "'  -3.42%  "
$ws.Range('E36').Value = "'  -1.75%  "
$ws.Range('E37').Value = "'  -2.52%  "
$ws.Range('E38').Value = "'  -2.62%  "
$ws.Range('E39').Value = "'  -2.42%  "
$ws.Range('D40').Value = "'1.93"
$ws.Range('E40').Value = "'  +2.44%  "
$ws.Range('D41').Value = "'0.510"
$ws.Range('E41').Value = "'  -2.47%  "
$ws.Range('E42').Value = "'  -0.14%  "
$ws.Range('E43').Value = "'  -2.24%  "
$ws.Range('D44').Value = "'0.0459"
$ws.Range('E44').Value = "'  +0.24%  "
$ws.Range('E45').Value = "'  -1.58%  "
$ws.Range('E46').Value = "'  -1.66%  "
$ws.Range('D47').Value = "'1.688.59"
$ws.Range('E47').Value = "'  -1.94%  "
$ws.Range('D48').Value = "'0.872"
$ws.Range('E48').Value = "'  -8.37%  "
$ws.Range('D49').Value = "'85.00"
$ws.Range('E49').Value = "'  -1.91%  "
$ws.Range('D50').Value = "'42.75"
$ws.Range('E50').Value = "'  +7.04%  "
$ws.Range('E51').Value = "'  -1.57%  "
